# Revert the Chart-of-Accounts import sample back to the richer
# "Accounts" template: rename the sheet, replace the header row, and
# replace all data rows with the new 11-column / 14-row dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Accounts"

# Every cell in this sheet is plain text (account codes like "10000" and
# the literal words TRUE/FALSE must stay text, not be coerced to numbers
# or booleans). A leading apostrophe forces Excel to store the cell as
# text while keeping the value itself exactly as written.
$header = @(
    "'AccountCode",
    "'Name",
    "'Type",
    "'Subtype",
    "'Description",
    "'ParentCode",
    "'IsSubledger",
    "'SubledgerType",
    "'FSLIBucket",
    "'InternalReportingBucket",
    "'Item"
)
for ($i = 0; $i -lt $header.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $header[$i]
}

$data = @(
    @("'10000", "'ASSETS", "'ASSET", "'Current", "'Assets Category", "'", "'FALSE", "'", "'Assets", "'Operations", "'"),
    @("'10100", "'Cash and Cash Equivalents", "'ASSET", "'Current", "'Cash and equivalents", "'10000", "'FALSE", "'", "'Assets", "'Operations", "'Cash"),
    @("'10101", "'Checking Account", "'ASSET", "'Current", "'Main checking account", "'10100", "'FALSE", "'", "'Assets", "'Operations", "'Banking"),
    @("'10102", "'Savings Account", "'ASSET", "'Current", "'Interest-bearing savings", "'10100", "'FALSE", "'", "'Assets", "'Operations", "'Banking"),
    @("'10200", "'Accounts Receivable", "'ASSET", "'Current", "'Money owed by customers", "'10000", "'TRUE", "'CUSTOMER", "'Assets", "'Sales", "'AR"),
    @("'20000", "'LIABILITIES", "'LIABILITY", "'Current", "'Liabilities Category", "'", "'FALSE", "'", "'Liabilities", "'Operations", "'"),
    @("'20100", "'Accounts Payable", "'LIABILITY", "'Current", "'Money owed to suppliers", "'20000", "'TRUE", "'VENDOR", "'Liabilities", "'Procurement", "'AP"),
    @("'30000", "'EQUITY", "'EQUITY", "'Retained Earnings", "'Equity Category", "'", "'FALSE", "'", "'Equity", "'Finance", "'"),
    @("'40000", "'REVENUE", "'REVENUE", "'Sales", "'Revenue Category", "'", "'FALSE", "'", "'Revenue", "'Sales", "'"),
    @("'40100", "'Services Revenue", "'REVENUE", "'Sales", "'Services income", "'40000", "'FALSE", "'", "'Revenue", "'Services", "'"),
    @("'50000", "'EXPENSES", "'EXPENSE", "'Operating", "'Expenses Category", "'", "'FALSE", "'", "'Expenses", "'Operations", "'"),
    @("'50100", "'Rent Expense", "'EXPENSE", "'Operating", "'Office rent costs", "'50000", "'FALSE", "'", "'Expenses", "'Facilities", "'Office"),
    @("'50200", "'Utilities Expense", "'EXPENSE", "'Operating", "'Utility costs", "'50000", "'FALSE", "'", "'Expenses", "'Facilities", "'Utilities"),
    @("'50300", "'Salaries Expense", "'EXPENSE", "'Operating", "'Employee salaries", "'50000", "'FALSE", "'", "'Expenses", "'HR", "'Payroll")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Clear any leftover cells outside the new A1:K15 range from the old
# A1:I11 layout (not applicable here since the new range is a strict
# superset, but kept defensively in case of future edits).
# Reset formatting to "Normal" so the quote-prefix used above to force
# text doesn't leave a stray cell style behind.
$ws.Range("A1:K15").Style = "Normal"
